$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror the value already present in A49 into the newly added A50 row,
# extending the data range by one row (as in the source diff).
$ws.Range("A50").Value = $ws.Range("A49").Value2
